$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(13)
$tr = $shape.TextFrame.TextRange
$tr.Characters(6, 6).Text = " 0 - 3"
